# Apply the CD-64 update to the lusid_holdings workbook:
#  - add a new "strategy" column (E) to the lusid_holdings sheet, populated
#    for every existing data row
#  - append three new holding rows (16-18) that also carry a strategy tag
#  - refresh the selection / active-sheet state left behind by the edit

$wb = $excel.ActiveWorkbook

$holdings = $wb.Worksheets.Item("lusid_holdings")
$instruments = $wb.Worksheets.Item("instruments")
$extHoldings = $wb.Worksheets.Item("ext_holdings")

# --- lusid_holdings: new "strategy" column -------------------------------
$holdings.Range("E1").Value = "strategy"

$strategyByRow = @{
    2  = "Quantitative"
    3  = "Rebalance"
    4  = "Quantitative"
    5  = "Quantitative"
    6  = "Rebalance"
    7  = "Quantitative"
    8  = "Growth"
    9  = "Quantitative"
    10 = "Growth"
    11 = "Rebalance"
    12 = "Quantitative"
    13 = "Quantitative"
    14 = "Quantitative"
    15 = "Quantitative"
}

foreach ($r in 2..15) {
    $holdings.Cells.Item($r, 5).Value = $strategyByRow[$r]
}

# --- lusid_holdings: three new holding rows -------------------------------
# Copy an existing data row first so the new rows inherit the same cell
# styling (the date column keeps its date number format) before the actual
# values are overwritten.
$newRows = @(
    @{ Row = 16; InstrumentId = "JE00B4T3BW64"; Units = 1450; Strategy = "Rebalance" }
    @{ Row = 17; InstrumentId = "GB0031743007"; Units = 790;  Strategy = "Rebalance" }
    @{ Row = 18; InstrumentId = "GB0005603997"; Units = 2300; Strategy = "Growth" }
)

foreach ($nr in $newRows) {
    $holdings.Range("A2:D2").Copy($holdings.Range("A$($nr.Row):D$($nr.Row)"))
    $holdings.Cells.Item($nr.Row, 2).Value = $nr.InstrumentId
    $holdings.Cells.Item($nr.Row, 3).Value = $nr.Units
    $holdings.Cells.Item($nr.Row, 4).Value = "GBP"
    $holdings.Cells.Item($nr.Row, 5).Value = $nr.Strategy
}

# Column D ("currency") no longer needs to be squeezed to a bestFit width
# now that column E exists next to it.
$holdings.Columns.Item(4).ColumnWidth = 10.33203125

# --- view state: selections on each sheet ---------------------------------
$instruments.Activate()
$instruments.Range("J13:J14").Select()

$extHoldings.Activate()
$extHoldings.Range("B13").Select()

$holdings.Activate()
$holdings.Range("E18").Select()
